# "break out stock.yaml completed"
#
# On the "day" sheet:
#   1. The four bsecode cells D841:D844 had been stored as text; re-enter
#      them as real numbers (same digits).
#   2. Append 19 freshly scraped rows (845-863), pushing the used range to
#      A1:I863. These rows keep bsecode (column D) as text, same as every
#      other bsecode cell on the sheet before this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. Fix up the bsecode values that used to be text -> numbers.
$ws.Range("D841").Value = 532541
$ws.Range("D842").Value = 533179
$ws.Range("D843").Value = 503806
$ws.Range("D844").Value = 531213

# --- 2. Append the new rows scraped on 06/11/2024.
$newRows = @(
    @(1,  "ULTRACEMCO", "Ultratech Cement Limited",                        "532538", 0.88,  11274.4,  282491,    "day", "06/11/2024 11:37:02"),
    @(2,  "BAJFINANCE",  "Bajaj Finance Limited",                          "500034", 1.09,  7006.2,   696268,    "day", "06/11/2024 11:37:02"),
    @(3,  "PERSISTENT",  "Persistent Systems Limited",                     "533179", 5.49,  5717.65,  690293,    "day", "06/11/2024 11:37:02"),
    @(4,  "PIIND",       "Pi Industries Limited",                          "523642", 1.68,  4614.2,   151790,    "day", "06/11/2024 11:37:02"),
    @(5,  "NAVINFLUOR",  "Navin Fluorine International Limited",           "532504", 3.56,  3495.8,   185776,    "day", "06/11/2024 11:37:02"),
    @(6,  "MPHASIS",     "Mphasis Limited",                                "526299", 1.35,  2896.4,   1687173,   "day", "06/11/2024 11:37:02"),
    @(7,  "SYNGENE",     "Syngene International Limited",                  "539268", 4.16,  903.85,   2578111,   "day", "06/11/2024 11:37:02"),
    @(8,  "RAMCOCEM",    "The Ramco Cements Limited",                      "500260", 0.77,  881.5,    254647,    "day", "06/11/2024 11:37:02"),
    @(9,  "CANFINHOME",  "Can Fin Homes Limited",                          "511196", 0.46,  869.8,    516437,    "day", "06/11/2024 11:37:02"),
    @(10, "ICICIPRULI",  "Icici Prudential Life Insurance Company Limited","540133", -2.17, 716.25,   864989,    "day", "06/11/2024 11:37:02"),
    @(11, "INDHOTEL",    "The Indian Hotels Company Limited",              "500850", 2.49,  684.65,   2440822,   "day", "06/11/2024 11:37:02"),
    @(12, "LICHSGFIN",   "Lic Housing Finance Limited",                    "500253", 0.37,  634.45,   956132,    "day", "06/11/2024 11:37:02"),
    @(13, "BALRAMCHIN",  "Balrampur Chini Mills Limited",                  "500038", -0.39, 607.95,   469331,    "day", "06/11/2024 11:37:02"),
    @(14, "RECLTD",      "Rural Electrification Corporation Limited",      "532955", 1.29,  532.8,    6153118,   "day", "06/11/2024 11:37:02"),
    @(15, "PFC",         "Power Finance Corporation Limited",              "532810", 1.31,  467.55,   7221536,   "day", "06/11/2024 11:37:02"),
    @(16, "M&MFIN",      "Mahindra & Mahindra Financial Services Limited", "532720", 0.18,  277.4,    638700,    "day", "06/11/2024 11:37:02"),
    @(17, "BANKBARODA",  "Bank Of Baroda",                                 "532134", 1.8,   262.55,   21915311,  "day", "06/11/2024 11:37:02"),
    @(18, "MANAPPURAM",  "Manappuram Finance Limited",                     "531213", -1.26, 156.98,   15482143,  "day", "06/11/2024 11:37:02"),
    @(19, "PNB",         "Punjab National Bank",                           "532461", 2.17,  106.98,   59984748,  "day", "06/11/2024 11:37:02")
)

$startRow = 845
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    # bsecode: keep as text even though it looks numeric (leading
    # apostrophe forces Excel to store it as a string, matching the
    # existing bsecode cells on this sheet).
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
